# EnterpriseJavaTimeLog.xlsx - "Cleans up comments in SearchExperiences servlet"
#
# Summary of changes applied:
#  1. D48 (Week 8 Project/Presentation note): reword the task description.
#  2. B48 (hours for that entry): 2.5 -> 8
#  3. D50: remove the leftover "2.5 doucmented plus Sun 10:10 - 11:45" note,
#     leaving the (styled) cell blank.
#  4. Move the sheet's scroll position / active selection up a few rows
#     (was topLeftCell B46 / selection D51, now B43 / D49).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1 & 2: update the task note text and the logged hours for row 48.
$ws.Range("D48").Value = "Project/Presentation: Got code working to backfill the indexes for Hibernate Search; completed several kinds of searchs; completed presentation materials"
$ws.Range("B48").Value = 8

# 3: clear out the old stray note in D50 (keeps its existing cell style).
$ws.Range("D50").ClearContents()

# 4: scroll the view up and move the active selection to D49.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("D49").Select()
